$d = $word.ActiveDocument

$replacements = @(
    @("32+6=", "5+37="),
    @("14+38=", "60+21="),
    @("61-23=", "96-48="),
    @("17+70=", "57-4="),
    @("69-68=", "67-25="),
    @("44-44=", "99-0="),
    @("60-18=", "65-65="),
    @("84-67=", "17+72="),
    @("92-49=", "14+77="),
    @("6+14=", "27+59="),
    @("36+42=", "30+14="),
    @("70-12=", "20+76="),
    @("28+61=", "1+88="),
    @("21+58=", "49-48="),
    @("96-80=", "19+37="),
    @("31+61=", "21-3="),
    @("59-15=", "30-4="),
    @("13+30=", "72-42="),
    @("51+6=", "49+19="),
    @("55-30=", "84-26="),
    @("32+58=", "48+28="),
    @("66-19=", "80-11="),
    @("63-33=", "63-8="),
    @("75-26=", "20+38="),
    @("35-13=", "43+6="),
    @("79-58=", "17+51="),
    @("48-24=", "51+10="),
    @("0+79=", "34-21="),
    @("93-17=", "48-22="),
    @("24+50=", "96-72="),
    @("9+46=", "63-3="),
    @("89-73=", "40+10="),
    @("7+63=", "67-33="),
    @("73+3=", "51+20="),
    @("97-75=", "77+1="),
    @("30+47=", "82-65="),
    @("50-38=", "23-18="),
    @("10+18=", "89-71="),
    @("61-27=", "70-64="),
    @("51-1=", "37+11="),
    @("0+49=", "39+26="),
    @("34-5=", "49+29="),
    @("92+7=", "55+39="),
    @("28+70=", "8+86="),
    @("70+29=", "59-32="),
    @("37-3=", "17+53="),
    @("46-20=", "80-41="),
    @("83-4=", "76-16="),
    @("16+62=", "8+79="),
    @("48+13=", "6+25="),
    @("63+3=", "0+67="),
    @("91-67=", "8+52="),
    @("39+9=", "21+1="),
    @("35-3=", "7+30="),
    @("62+35=", "10+26="),
    @("32+12=", "9+31="),
    @("20+4=", "82-7="),
    @("14+23=", "31+31="),
    @("14-2=", "62+28="),
    @("2+3=", "27+39="),
    @("93-40=", "52-47="),
    @("56-0=", "86-64="),
    @("59-25=", "14+43="),
    @("51+11=", "6+61="),
    @("47+17=", "12+78="),
    @("55-17=", "40-19="),
    @("11+58=", "63+27="),
    @("56-32=", "21+64="),
    @("22-20=", "19+31="),
    @("15+57=", "82-27="),
    @("96-69=", "20+70="),
    @("87-6=", "16+48="),
    @("4+55=", "63+11="),
    @("30-25=", "80-60="),
    @("93-90=", "54-30="),
    @("49+36=", "93-5="),
    @("52-10=", "63-53="),
    @("20+79=", "75-21="),
    @("57-16=", "52+26="),
    @("80-38=", "3+3="),
    @("33+18=", "68-42="),
    @("92-63=", "63-20="),
    @("93-32=", "37-14="),
    @("98-84=", "37+30="),
    @("73+21=", "13-9="),
    @("27+19=", "12+5="),
    @("66+32=", "66+21="),
    @("32+63=", "10+88="),
    @("93+2=", "4+88="),
    @("36+52=", "74-30="),
    @("60-11=", "75-27="),
    @("62-37=", "58-41="),
    @("74-38=", "67-42="),
    @("3+93=", "58-30="),
    @("41+7=", "26+51="),
    @("65-46=", "34+61="),
    @("45-35=", "85-54="),
    @("8-1=", "3+55="),
    @("58+9=", "88-75="),
    @("28+62=", "56+3="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()